# "Small Fixes after Merging Sprint 2 and Sprint 3"
# The TODAY() formulas that had been filled down E4:E35 are converted to
# plain static date values (E2:E3 are left as live "=TODAY()" formulas).
# The sheet selection/scroll position is also updated to reflect the range
# that was just edited (E4:E35, scrolled so row 9 is at the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Freeze the filled-down TODAY() formulas in E4:E35 into plain values
# (paste-special-values over the formula range), matching the diff where
# the <f> elements disappear and only <v>44422</v> remains for rows 4-35.
$ws.Range("E4:E35").Value = 44422

# Reflect the post-edit selection/scroll state: the user had E4:E35
# selected (active cell E4) with the view scrolled down so row 9 is the
# first visible row.
$ws.Range("E4:E35").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
